# Daily attendance processing - 2025-10-07 01:12:57
# Reorder the "Recorded By" (column G) entries so that "System" always
# appears first in the comma-separated list, preserving the relative
# order of the remaining entries. The comparison must be case sensitive,
# since "system" (lowercase) and "System" are distinct, separate entries.

function Test-CaseSensitiveEquals($strA, $strB) {
    if ($strA.Length -ne $strB.Length) { return $false }
    for ($ci = 0; $ci -lt $strA.Length; $ci++) {
        $codeA = [int][char]$strA[$ci]
        $codeB = [int][char]$strB[$ci]
        if ($codeA -ne $codeB) { return $false }
    }
    return $true
}

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)  # column G
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -eq "") { continue }

    $parts = @($val -split ",\s*")

    $systemIndex = -1
    for ($pi = 0; $pi -lt $parts.Count; $pi++) {
        if (Test-CaseSensitiveEquals $parts[$pi] "System") {
            $systemIndex = $pi
            break
        }
    }

    if ($systemIndex -ge 0 -and $systemIndex -ne 0) {
        $rest = New-Object System.Collections.ArrayList
        for ($qi = 0; $qi -lt $parts.Count; $qi++) {
            if ($qi -ne $systemIndex) {
                [void]$rest.Add($parts[$qi])
            }
        }
        $newParts = @("System") + $rest
        $cell.Value2 = [string]::Join(", ", $newParts)
    }
}
